$d = $word.ActiveDocument

$replacements = @(
    @("257÷6=", "826÷6="),
    @("571÷2=", "933÷9="),
    @("751÷3=", "128÷3="),
    @("858÷4=", "264÷4="),
    @("207÷8=", "738÷8="),
    @("667÷2=", "865÷8="),
    @("791÷8=", "959÷9="),
    @("535÷4=", "743÷6="),
    @("604÷2=", "974÷8="),
    @("688÷2=", "224÷3="),
    @("289÷2=", "540÷3="),
    @("585÷8=", "929÷5="),
    @("278÷5=", "955÷9="),
    @("928÷5=", "133÷4="),
    @("811÷3=", "797÷7="),
    @("268÷5=", "736÷6="),
    @("297÷8=", "331÷7="),
    @("486÷8=", "383÷9="),
    @("595÷8=", "973÷5="),
    @("382÷5=", "286÷9="),
    @("693÷8=", "281÷6="),
    @("970÷8=", "240÷7="),
    @("428÷6=", "511÷3="),
    @("850÷6=", "836÷2="),
    @("248÷6=", "530÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
